# Swap the "category" and "group" columns in the SectorGroup codelist sheet.
#
# Before the edit, column D held codeforiati:category-name / category values
# and column E held codeforiati:group-name / group values (and similarly F
# held group-code while G held category-code). The commit swaps these so
# that D/E (and F/G) exchange places for every row, including the header
# row, so that D becomes "group" and E becomes "category".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $dCell = $ws.Cells.Item($r, 4)
    $eCell = $ws.Cells.Item($r, 5)
    $fCell = $ws.Cells.Item($r, 6)
    $gCell = $ws.Cells.Item($r, 7)

    $dVal = $dCell.Value()
    $eVal = $eCell.Value()
    $fVal = $fCell.Value()
    $gVal = $gCell.Value()

    $dCell.Value = $eVal
    $eCell.Value = $dVal
    $fCell.Value = $gVal
    $gCell.Value = $fVal
}
